# "queries de tudo ou nada - categoria-mencao"
# Re-orders the UF rows that tie on a zero/empty "all-or-nothing" metric.
# The underlying query result order for these tied UFs changed, which (because
# PA/ES and MT/TO also swap table position in the shared-strings table as a
# side effect) also flips the PA/ES and MT/TO rows on the "uf-qtd" sheet even
# though their numeric values stay put.

$wb = $excel.ActiveWorkbook

# --- uf-qtd ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("uf-qtd")
$ws.Range("A11").Value = "PA"
$ws.Range("A12").Value = "ES"
$ws.Range("A26").Value = "MT"
$ws.Range("A27").Value = "TO"

# --- uf-tot-arrecad ---------------------------------------------------------
$ws = $wb.Worksheets.Item("uf-tot-arrecad")
$ws.Range("A21").Value = "RO"
$ws.Range("A22").Value = "AP"
$ws.Range("A23").Value = "TO"
$ws.Range("A24").Value = "AL"
$ws.Range("A25").Value = "MT"
$ws.Range("A26").Value = "PI"
$ws.Range("A27").Value = "RN"

# --- uf-avg-arrecad -----------------------------------------------------
$ws = $wb.Worksheets.Item("uf-avg-arrecad")
$ws.Range("A21").Value = "AL"
$ws.Range("A22").Value = "MT"
$ws.Range("A23").Value = "RO"
$ws.Range("A24").Value = "AP"
$ws.Range("A25").Value = "TO"
$ws.Range("A26").Value = "PI"
$ws.Range("A27").Value = "RN"

# --- uf-max-arrecad -----------------------------------------------------
$ws = $wb.Worksheets.Item("uf-max-arrecad")
$ws.Range("A21").Value = "AL"
$ws.Range("A22").Value = "MT"
$ws.Range("A23").Value = "RO"
$ws.Range("A24").Value = "AP"
$ws.Range("A25").Value = "TO"
$ws.Range("A26").Value = "PI"
$ws.Range("A27").Value = "RN"

# --- uf-tx-sucesso ----------------------------------------------------------
$ws = $wb.Worksheets.Item("uf-tx-sucesso")
$ws.Range("A21").Value = "RO"
$ws.Range("A22").Value = "AP"
$ws.Range("A23").Value = "TO"
$ws.Range("A24").Value = "AL"
$ws.Range("A25").Value = "MT"
$ws.Range("A26").Value = "PI"
$ws.Range("A27").Value = "RN"
